$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.142.85"
$ws.Range("E2").Value = "  +2.19%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.531.15"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.50"
$ws.Range("E5").Value = "  +1.52%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.44"
$ws.Range("E6").Value = "  +3.52%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +0.95%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.531.02"
$ws.Range("E9").Value = "  +0.39%  "

# Row 10
$ws.Range("E10").Value = "  +2.34%  "

# Row 11
$ws.Range("E11").Value = "  +2.59%  "

# Row 12
$ws.Range("E12").Value = "  -0.32%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.346"
$ws.Range("E13").Value = "  -1.17%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.85"
$ws.Range("E14").Value = "  +0.33%  "

# Row 15
$ws.Range("E15").Value = "  +2.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.990.46"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.297.61"
$ws.Range("E17").Value = "  +2.76%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.528.40"
$ws.Range("E18").Value = "  +0.03%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.98"
$ws.Range("E19").Value = "  +1.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.54"
$ws.Range("E20").Value = "  +1.77%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "367.13"
$ws.Range("E21").Value = "  +5.51%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  -0.18%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.71"
$ws.Range("E23").Value = "  +0.99%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.07%  "

# Row 25
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("E25").Value = "  -2.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.91"
$ws.Range("E26").Value = "  +0.98%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.18"
$ws.Range("E27").Value = "  +2.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.661.26"
$ws.Range("E28").Value = "  +1.02%  "

# Row 29
$ws.Range("E29").Value = "  -0.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0996"
$ws.Range("E30").Value = "  +1.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "540.75"
$ws.Range("E31").Value = "  +2.65%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.28"
$ws.Range("E32").Value = "  +1.62%  "

# Row 33
$ws.Range("E33").Value = "  +1.07%  "

# Row 34
$ws.Range("E34").Value = "  +2.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.129"
$ws.Range("E35").Value = "  -1.14%  "

# Row 36
$ws.Range("E36").Value = "  +0.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.46"
$ws.Range("E37").Value = "  -0.26%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.09"
$ws.Range("E38").Value = "  -0.28%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.85"
$ws.Range("E39").Value = "  +1.11%  "

# Row 40
$ws.Range("E40").Value = "  +1.65%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.355"
$ws.Range("E41").Value = "  +0.03%  "

# Row 42
$ws.Range("E42").Value = "  +0.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.21"
$ws.Range("E43").Value = "  +2.11%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.55"
$ws.Range("E44").Value = "  +1.74%  "

# Row 45
$ws.Range("E45").Value = "  -0.02%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.21"
$ws.Range("E46").Value = "  -1.59%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.560"
$ws.Range("E47").Value = "  +0.08%  "

# Row 48
$ws.Range("E48").Value = "  +1.00%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0277"
$ws.Range("E49").Value = "  +3.14%  "

# Row 50
$ws.Range("E50").Value = "  -1.15%  "

# Row 51
$ws.Range("E51").Value = "  -0.40%  "
